$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet set-up: rename Sheet1 -> testSheet, add testSheet2 right after it
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "testSheet"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "testSheet2"

# =========================================================================
# testSheet ("site / name / status / date / comment" form submissions)
# =========================================================================

$headers1 = @("site", "name", "status", "date", "comment")
for ($c = 1; $c -le 5; $c++) {
    $ws1.Cells.Item(1, $c).Value = $headers1[$c - 1]
}

$row1 = @("Downtown", "Bryce Eppler", "status goes here", "12/06/2022:15:64:12PST", "Today was a good day for all the guys at work :)")
for ($r = 2; $r -le 5; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws1.Cells.Item($r, $c).Value = $row1[$c - 1]
    }
}

# Column width (applies to the whole column, matching the sheet-wide default)
$ws1.Columns.Item(1).ColumnWidth = 34.998697916666664

# Row heights for the used rows
$ws1.Range("A1:E5").RowHeight = 30

# Header row: bold font, bottom border, wrapped text
$headerRange1 = $ws1.Range("A1:E1")
$headerRange1.Font.Bold = $true
$headerRange1.WrapText = $true
$headerRange1.Borders.Item(9).LineStyle = 1

# Data rows: regular font, wrapped text
$dataRange1 = $ws1.Range("A2:E5")
$dataRange1.WrapText = $true

$ws1.Range("B9").Select()

# =========================================================================
# testSheet2 (ID / Site / Name / Status / Tasks Completed / Comments / Date)
# =========================================================================

$headers2 = @("ID", "Site", "Name", "Status", "Tasks Completed", "Comments", "Date")
for ($c = 1; $c -le 7; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers2[$c - 1]
}

$headerRange2 = $ws2.Range("A1:G1")
$headerRange2.Font.Bold = $true

$ws2.Range("G1").HorizontalAlignment = -4152

$ws2.Range("A2").Value = 1

$ws2.Range("A3").Value = "Downtown"
$ws2.Range("B3").Value = "Bryce Eppler"
$ws2.Range("C3").Value = "status goes here"
$ws2.Range("D3").Value = "12/06/2022:15:64:12PST"
$ws2.Range("E3").Value = "Today was a good day for all the guys at work :)"

$ws2.Columns.Item(5).ColumnWidth = 14.498697916666664

$ws2.Range("G2").Select()
